# "update my journal de board"
# Edits the "Feuil1" journal table: wrap two long remarks onto a second
# line (Alt+Enter inside the cell), tone down one remark's wording, and
# nudge column A's width back down.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets("Feuil1")

# --- C10: soften "un petit erreur" -> "un erreur" ----------------------
$ws.Range("C10").Value = "un erreur au niveau gestion de temps  "

# --- B3: break the "Reunion avec l'expert ..." note onto two lines ----
$ws.Range("B3").Value = "Reunion avec l'expert + signature du cahier des charges de la part du `ncandidat durant la réunion"
$ws.Range("B3").WrapText = $true

# --- B13: break the "faire une modification ..." note onto two lines --
$ws.Range("B13").Value = "faire une modification sur la base de `ndonnées un rendez-vous avec le chef de projet "
$ws.Range("B13").WrapText = $true
$ws.Rows.Item(13).RowHeight = 31.5

# --- Column A: slightly narrower now that the table reflows -----------
$ws.Columns.Item(1).ColumnWidth = 9.8
